$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 98; $row++) {
    $ws.Range("K$row").Formula = "=E$row/D$row"
    $ws.Range("L$row").Formula = "=H$row/F$row"
}

$null = $ws.Range("P6").Select()
